# Update the handback status timestamps on the zh-cn and de-de sheets
# (regenerated report run produced new Correspond Handoff / Handback timestamps).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-20 16:39:42"
$zhcn.Range("H3").Value = "2016-03-20 16:40:04"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-20 16:39:45"
$dede.Range("H3").Value = "2016-03-20 16:40:10"
